# Updated symbol list on Thu Feb  2 19:58:30 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Price/Volume cells store plain text (e.g. "328.83", "5.38%"), not numbers.
# Force text format before writing so Excel does not auto-convert the value into
# a number/percentage, then restore the default "Normal" style so no stray cell
# formatting (s=) is left behind -- only the text content should change.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "328.83"
Set-TextValue "E2" "5.38%"

Set-TextValue "D3" "40.47"
Set-TextValue "E3" "9.86%"

Set-TextValue "D4" "5.833"
Set-TextValue "E4" "14.04%"

Set-TextValue "D5" "0.08135"
Set-TextValue "E5" "4.24%"

Set-TextValue "D6" "4.602"
Set-TextValue "E6" "5.03%"

Set-TextValue "D7" "8.756"
Set-TextValue "E7" "4.59%"

Set-TextValue "D8" "1.969"
Set-TextValue "E8" "5.30%"

Set-TextValue "E9" "-0.28%"

Set-TextValue "D10" "0.9461"
Set-TextValue "E10" "1.75%"

Set-TextValue "D11" "0.1324"
Set-TextValue "E11" "11.94%"

Set-TextValue "D12" "0.1997"
Set-TextValue "E12" "4.85%"

Set-TextValue "D13" "8.969"
Set-TextValue "E13" "42.53%"

Set-TextValue "D14" "0.09345"
Set-TextValue "E14" "5.32%"

Set-TextValue "D15" "0.03451"
Set-TextValue "E15" "5.34%"

Set-TextValue "D16" "0.09601"
Set-TextValue "E16" "0.41%"

Set-TextValue "D17" "0.001324"
Set-TextValue "E17" "-3.56%"

Set-TextValue "D18" "0.006304"
Set-TextValue "E18" "6.18%"

Set-TextValue "D19" "3.364"
Set-TextValue "E19" "0.53%"

Set-TextValue "D20" "0.3541"
Set-TextValue "E20" "2.76%"

Set-TextValue "D21" "0.1402"
Set-TextValue "E21" "8.78%"

Set-TextValue "D22" "0.2412"
Set-TextValue "E22" "5.20%"

Set-TextValue "D23" "0.04437"
Set-TextValue "E23" "2.81%"

Set-TextValue "E24" "5.94%"

Set-TextValue "D25" "0.004423"
Set-TextValue "E25" "1.82%"

Set-TextValue "D26" "0.0001093"
Set-TextValue "E26" "-17.15%"

Set-TextValue "D27" "0.0003995"
Set-TextValue "E27" "1.13%"

Set-TextValue "D39" "0.02471"
Set-TextValue "E39" "10.80%"

Set-TextValue "D40" "0.05292"
Set-TextValue "E40" "3.55%"

Set-TextValue "D41" "0.007473"
Set-TextValue "E41" "-1.96%"

Set-TextValue "D42" "0.1438"
Set-TextValue "E42" "4.45%"

Set-TextValue "D43" "0.009019"
Set-TextValue "E43" "9.81%"

Set-TextValue "D44" "0.002056"
Set-TextValue "E44" "3.75%"

Set-TextValue "D45" "0.01053"
Set-TextValue "E45" "35.61%"

Set-TextValue "E46" "8.49%"

Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "1.06%"

Set-TextValue "D48" "0.003509"
Set-TextValue "E48" "23.34%"

Set-TextValue "D49" "0.001803"
Set-TextValue "E49" "7.49%"

Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "1.06%"

Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "1.06%"
